# Re-pull data / push all data / mean calculation
# Update column F ("dSF") values for the rows whose underlying data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -10
    3  = 7
    6  = -2
    9  = -2
    10 = -4
    11 = 3
    13 = -2
    22 = -4
    23 = 0
    26 = 6
    27 = -1
    28 = 2
    29 = -6
    33 = -2
    35 = -2
    36 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
